$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q4" and push the
#    existing 2021-Q2 / 2021-Q1 / 2020-Q4 rows down by one.
# -----------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing data rows down (bottom-up to avoid clobbering).
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# Fill the new row 2 with the 2022-Q4 summary values.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.18

# Re-number column A (the 0-based index column) for the shifted rows.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# -----------------------------------------------------------------
# 2) Add the new "2022-Q4" sheet, positioned right after "总计".
#    Clone the "2021-Q2" sheet so we inherit identical sheet-level
#    formatting (sheetPr/pageMargins/styles) instead of building a
#    blank sheet from scratch.
# -----------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q2")
$template.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Header row tweak: this quarter's sheet calls column D "基金规模"
# instead of "基金金额".
$q4.Range("D1").Value = "基金规模"

# The template only has 4 data rows (rows 2-5); we need 7 (rows 2-8).
# Clone row 5's formatting down to create rows 6-8.
$q4.Range("A5:H5").Copy($q4.Range("A6:H6"))
$q4.Range("A5:H5").Copy($q4.Range("A7:H7"))
$q4.Range("A5:H5").Copy($q4.Range("A8:H8"))

# Columns B, D, E, F, G hold numeric-looking text (fund codes with
# leading zeros, percentages formatted as strings, etc.) in this
# workbook's convention, so force text formatting before writing
# them to avoid silent numeric coercion. (Union ranges such as
# "B2:B8,D2:G8" aren't reliably honoured by NumberFormat, so each
# contiguous block is set individually.)
$q4.Range("B2:B8").NumberFormat = "@"
$q4.Range("D2:G8").NumberFormat = "@"

$data = @(
    @(0, "014155", "国泰君安中证500指数增强A", "7.70", "92.93", "1.03", "0.0793", 9),
    @(1, "014156", "国泰君安中证500指数增强C", "4.81", "92.93", "1.03", "0.0495", 9),
    @(2, "006441", "中信建投中证500指数增强C", "2.27", "93.50", "0.83", "0.0188", 6),
    @(3, "006440", "中信建投中证500指数增强A", "2.19", "93.50", "0.83", "0.0182", 6),
    @(4, "159990", "银华巨潮小盘价值ETF", "0.78", "97.02", "1.28", "0.0100", 6),
    @(5, "003717", "中银量化精选灵活配置混合A", "0.41", "90.94", "1.20", "0.0049", 9),
    @(6, "010484", "中银量化精选灵活配置混合C", "0.02", "90.94", "1.20", "0.0002", 9)
)

$r = 2
foreach ($row in $data) {
    $q4.Range("A$r").Value = $row[0]
    $q4.Range("B$r").Value = $row[1]
    $q4.Range("C$r").Value = $row[2]
    $q4.Range("D$r").Value = $row[3]
    $q4.Range("E$r").Value = $row[4]
    $q4.Range("F$r").Value = $row[5]
    $q4.Range("G$r").Value = $row[6]
    $q4.Range("H$r").Value = $row[7]
    $r = $r + 1
}
